$wb = $excel.ActiveWorkbook

# The workbook's originally saved active tab is sheet index 8 (1-based),
# i.e. "TGZ-S-48-100_300RI". Selecting a range on another sheet below
# switches the active sheet/tab, so we re-activate this one at the end
# to keep the workbook's active tab unchanged.
$originalActiveSheetIndex = 8

# Rename sheets: "xyzRI" -> "xyz-RI"
$renames = @{
    "TGZ-S-48-50_100RI"  = "TGZ-S-48-50_100-RI"
    "TGZ-S-48-100_250RI" = "TGZ-S-48-100_250-RI"
    "TGZ-S-48-100_300RI" = "TGZ-S-48-100_300-RI"
}

foreach ($oldName in $renames.Keys) {
    $newName = $renames[$oldName]
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $newName
}

# Update the selected cell on the renamed "TGZ-S-48-50_100-RI" sheet from A11 to G26
$ws = $wb.Worksheets.Item("TGZ-S-48-50_100-RI")
$ws.Range("G26").Select()

# Restore the originally active sheet/tab (selecting above switched it)
$origSheet = $wb.Worksheets.Item($originalActiveSheetIndex)
$origSheet.Activate()
